$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Text Mining Grundlagen"
$ws.Range("B5").Value = "Fallstudie Populismus"
$ws.Range("B6").Value = "Word Embeddings"
$ws.Range("D6").Value = "Am 3. November entfällt der Unterricht in diesem Modul."
$ws.Range("B7").Value = "Projektwoche Twitter Hate Speech"
$ws.Range("B9").Value = "Hate Speech - Stand der Forschung"
$ws.Range("B10").Value = "Regression"
$ws.Range("B11").Value = "Klassifikation"
$ws.Range("B12").Value = "Projektwoche Twitter Hate Speech 2"
$ws.Range("B13").Value = "Quarto Blog"

$ws.Range("B9").Select()
